# Add team record (Wins/Losses/Ties) columns to the LAD_1996 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells in row 1
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Match the header formatting used by the rest of row 1 (bold, centered, bordered)
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the team record for every data row (2 through 42)
$lastRow = 42
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 29).Value = 90
    $ws.Cells.Item($r, 30).Value = 72
    $ws.Cells.Item($r, 31).Value = 0
}
